# DatasetTable.xlsx edit: restructure the spatial/temporal resolution columns
# - Swap TempRes/SpatCov column order (D/E)
# - Insert new "SpatRes" column (F), derived from the spatial-resolution portion
#   of the old free-text "Text" column (now removed)
# - Shift the remaining boolean coverage columns (SurfTemp..Snow) one column right
# - Update row 8's TempRes value from "One day" to "One day each month"
# - Resize column D, drop the custom width that used to belong to column O
# - Update the active-cell selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = New-Object 'object[,]' 10,12

# Row 1 (header)
$values[0,0]  = "TempRes"
$values[0,1]  = "SpatCov"
$values[0,2]  = "SpatRes"
$values[0,3]  = "AirTemp"
$values[0,4]  = "SurfTemp"
$values[0,5]  = "SoilTemp"
$values[0,6]  = "Radiation"
$values[0,7]  = "Wind"
$values[0,8]  = "Precipitation"
$values[0,9]  = "Humidity"
$values[0,10] = "SoilMoist"
$values[0,11] = "Snow"

# Row 2 - SCAN
$values[1,0]  = "Hourly"
$values[1,1]  = "US"
$values[1,2]  = "Stations"
$values[1,3]  = "T"
$values[1,4]  = "F"
$values[1,5]  = "T"
$values[1,6]  = "T"
$values[1,7]  = "T"
$values[1,8]  = "T"
$values[1,9]  = "T"
$values[1,10] = "T"
$values[1,11] = "F"

# Row 3 - ERA5
$values[2,0]  = "Hourly"
$values[2,1]  = "Global"
$values[2,2]  = "0.1° x 0.1°"
$values[2,3]  = "T"
$values[2,4]  = "T"
$values[2,5]  = "T"
$values[2,6]  = "T"
$values[2,7]  = "T"
$values[2,8]  = "T"
$values[2,9]  = "F"
$values[2,10] = "F"
$values[2,11] = "T"

# Row 4 - GLDAS
$values[3,0]  = "3-hourly"
$values[3,1]  = "Global"
$values[3,2]  = "0.25° x 0.25°"
$values[3,3]  = "T"
$values[3,4]  = "T"
$values[3,5]  = "T"
$values[3,6]  = "T"
$values[3,7]  = "T"
$values[3,8]  = "T"
$values[3,9]  = "T"
$values[3,10] = "T"
$values[3,11] = "T"

# Row 5 - GRIDMET
$values[4,0]  = "Daily"
$values[4,1]  = "US"
$values[4,2]  = "Stations"
$values[4,3]  = "T"
$values[4,4]  = "F"
$values[4,5]  = "F"
$values[4,6]  = "T"
$values[4,7]  = "T"
$values[4,8]  = "T"
$values[4,9]  = "F"
$values[4,10] = "F"
$values[4,11] = "F"

# Row 6 - NOAA_NCDC
$values[5,0]  = "Daily"
$values[5,1]  = "Global"
$values[5,2]  = "Stations"
$values[5,3]  = "T"
$values[5,4]  = "F"
$values[5,5]  = "F"
$values[5,6]  = "F"
$values[5,7]  = "F"
$values[5,8]  = "T"
$values[5,9]  = "F"
$values[5,10] = "F"
$values[5,11] = "T"

# Row 7 - microclimUS
$values[6,0]  = "Hourly"
$values[6,1]  = "US"
$values[6,2]  = "0.6° x 0.6°"
$values[6,3]  = "T"
$values[6,4]  = "T"
$values[6,5]  = "T"
$values[6,6]  = "T"
$values[6,7]  = "F"
$values[6,8]  = "F"
$values[6,9]  = "T"
$values[6,10] = "T"
$values[6,11] = "F"

# Row 8 - microclim
$values[7,0]  = "One day each month"
$values[7,1]  = "Global"
$values[7,2]  = "0.17° x 0.17°"
$values[7,3]  = "T"
$values[7,4]  = "T"
$values[7,5]  = "T"
$values[7,6]  = "T"
$values[7,7]  = "T"
$values[7,8]  = "F"
$values[7,9]  = "T"
$values[7,10] = "F"
$values[7,11] = "F"

# Row 9 - SNODAS
$values[8,0]  = "Daily"
$values[8,1]  = "US"
$values[8,2]  = "0.01° x 0.01°"
$values[8,3]  = "F"
$values[8,4]  = "F"
$values[8,5]  = "F"
$values[8,6]  = "F"
$values[8,7]  = "F"
$values[8,8]  = "F"
$values[8,9]  = "F"
$values[8,10] = "F"
$values[8,11] = "T"

# Row 10 - NicheMapR
$values[9,0]  = "Hourly"
$values[9,1]  = "Global"
$values[9,2]  = "30m?"
$values[9,3]  = "T"
$values[9,4]  = "T"
$values[9,5]  = "T"
$values[9,6]  = "T"
$values[9,7]  = "T"
$values[9,8]  = "F"
$values[9,9]  = "T"
$values[9,10] = "F"
$values[9,11] = "T"

$ws.Range("D1:O10").Value = $values

# Column D now holds the TempRes values and is sized to fit; column O no
# longer needs the wide custom width that used to hold the free-text column.
$ws.Columns.Item(4).ColumnWidth = 17.1
$ws.Columns.Item(15).ColumnWidth = 8.43

# Update the saved selection to match the edited workbook.
$ws.Range("D10").Select() | Out-Null
